# Adapt column header formatting to respective input file names (#7):
#   "<header>_old" -> "<header>_FV2404"
#   "<header>_new" -> "<header>_FV2410"
# Also turn the header+data range into a real Excel Table ("Table1")
# and freeze the header row, matching the upstream regenerated export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

# Columns A-J carry the "_old" suffix -> rename to "_FV2404"
$col = 1
foreach ($h in $headers) {
    $ws.Cells.Item(1, $col).Value = "$h" + "_FV2404"
    $col = $col + 1
}

# Column K ("diff") is unaffected.

# Columns L-U carry the "_new" suffix -> rename to "_FV2410"
$col = 12
foreach ($h in $headers) {
    $ws.Cells.Item(1, $col).Value = "$h" + "_FV2410"
    $col = $col + 1
}

# Convert the whole used range (header row + 56 data rows, A:U) into an
# Excel Table named "Table1" with a plain style (no accent color).
$range = $ws.Range("A1:U57")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# Freeze the header row (split below row 1) like the regenerated sheet.
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
